# Updated symbol list on Sun Dec 25 14:30:58 UTC 2022 with GitHub Actions
#
# Applies the latest price/volume-label refresh to the crypto tracker sheet.
# Price values in column D are stored as literal text (not numbers) in the
# source data, so each new value is written with a leading apostrophe to
# force text entry, then the cell style is reset to "Normal" so Excel's
# "number stored as text" quote-prefix marker doesn't leave a stray style
# behind on cells that previously had none.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $ws.Range($cellRef).Value = "'" + $newValue
    $ws.Range($cellRef).Style = "Normal"
}

# Column D ("Price") updates
Set-TextValue "D2"  "244.03"
Set-TextValue "D4"  "5.414"
Set-TextValue "D5"  "0.05945"
Set-TextValue "D6"  "3.394"
Set-TextValue "D8"  "0.9265"
Set-TextValue "D9"  "0.1417"
Set-TextValue "D10" "0.07423"
Set-TextValue "D12" "0.03077"
Set-TextValue "D13" "0.09349"
Set-TextValue "D15" "0.001596"
Set-TextValue "D16" "0.04823"
Set-TextValue "D17" "0.0005943"
Set-TextValue "D18" "0.005501"
Set-TextValue "D19" "0.004328"
Set-TextValue "D20" "0.0009833"
Set-TextValue "D21" "0.00007504"
Set-TextValue "D22" "3.659"
Set-TextValue "D23" "6.452"
Set-TextValue "D24" "2.185"
Set-TextValue "D25" "0.3248"
Set-TextValue "D26" "0.1340"
Set-TextValue "D27" "0.0002447"
Set-TextValue "D40" "0.03905"
Set-TextValue "D41" "0.006214"
Set-TextValue "D42" "0.1072"
Set-TextValue "D44" "0.007284"
Set-TextValue "D45" "0.00005185"
Set-TextValue "D47" "0.0005803"
Set-TextValue "D49" "0.002317"

# Column E ("Volume(1h)") label updates
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"
